{"js": "// Update the date line and every \"dividend\u00f7divisor=\" expression in the\n// three-digit \u00f7 one-digit practice sheet to the new day's values.\nconst replacements = [\n  [\"2024-09-04 Wednesday\", \"2024-09-05 Thursday\"],\n  [\"638\u00f78=\", \"262\u00f76=\"],\n  [\"369\u00f74=\", \"340\u00f73=\"],\n  [\"747\u00f72=\", \"963\u00f73=\"],\n  [\"582\u00f72=\", \"569\u00f79=\"],\n  [\"561\u00f75=\", \"745\u00f76=\"],\n  [\"505\u00f76=\", \"236\u00f78=\"],\n  [\"544\u00f75=\", \"637\u00f74=\"],\n  [\"331\u00f72=\", \"709\u00f77=\"],\n  [\"793\u00f72=\", \"144\u00f72=\"],\n  [\"604\u00f75=\", \"995\u00f74=\"],\n  [\"623\u00f73=\", \"752\u00f74=\"],\n  [\"604\u00f79=\", \"620\u00f73=\"],\n  [\"849\u00f79=\", \"168\u00f74=\"],\n  [\"631\u00f74=\", \"515\u00f79=\"],\n  [\"984\u00f76=\", \"424\u00f79=\"],\n  [\"425\u00f77=\", \"254\u00f75=\"],\n  [\"543\u00f75=\", \"363\u00f77=\"],\n  [\"385\u00f77=\", \"130\u00f77=\"],\n  [\"421\u00f72=\", \"604\u00f76=\"],\n  [\"754\u00f76=\", \"616\u00f79=\"],\n  [\"265\u00f79=\", \"444\u00f72=\"],\n  [\"509\u00f73=\", \"819\u00f73=\"],\n  [\"486\u00f78=\", \"956\u00f77=\"],\n  [\"301\u00f74=\", \"581\u00f74=\"],\n  [\"910\u00f72=\", \"809\u00f77=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  for (const r of found.items) {\n    r.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date line and every \"dividend\u00f7divisor=\" expression in the\n# three-digit \u00f7 one-digit practice sheet to the new day's values.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-09-04 Wednesday\", \"2024-09-05 Thursday\"),\n    @(\"638\u00f78=\", \"262\u00f76=\"),\n    @(\"369\u00f74=\", \"340\u00f73=\"),\n    @(\"747\u00f72=\", \"963\u00f73=\"),\n    @(\"582\u00f72=\", \"569\u00f79=\"),\n    @(\"561\u00f75=\", \"745\u00f76=\"),\n    @(\"505\u00f76=\", \"236\u00f78=\"),\n    @(\"544\u00f75=\", \"637\u00f74=\"),\n    @(\"331\u00f72=\", \"709\u00f77=\"),\n    @(\"793\u00f72=\", \"144\u00f72=\"),\n    @(\"604\u00f75=\", \"995\u00f74=\"),\n    @(\"623\u00f73=\", \"752\u00f74=\"),\n    @(\"604\u00f79=\", \"620\u00f73=\"),\n    @(\"849\u00f79=\", \"168\u00f74=\"),\n    @(\"631\u00f74=\", \"515\u00f79=\"),\n    @(\"984\u00f76=\", \"424\u00f79=\"),\n    @(\"425\u00f77=\", \"254\u00f75=\"),\n    @(\"543\u00f75=\", \"363\u00f77=\"),\n    @(\"385\u00f77=\", \"130\u00f77=\"),\n    @(\"421\u00f72=\", \"604\u00f76=\"),\n    @(\"754\u00f76=\", \"616\u00f79=\"),\n    @(\"265\u00f79=\", \"444\u00f72=\"),\n    @(\"509\u00f73=\", \"819\u00f73=\"),\n    @(\"486\u00f78=\", \"956\u00f77=\"),\n    @(\"301\u00f74=\", \"581\u00f74=\"),\n    @(\"910\u00f72=\", \"809\u00f77=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n\n$d.Saved = $false\n"}
